$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = 0.83379886465449327
$ws.Range("AH1").Value = 0.77369751122668884
$ws.Range("BN1").Value = 0.56077363481638542
$ws.Range("BO1").Value = 0.76926545454947814
$ws.Range("BP1").Value = 0.93812034679608414
$ws.Range("A2").Value = 0.97375977933374092
$ws.Range("C2").Value = 0.83268421090543077
$ws.Range("D2").Value = 0.67197975706414836
$ws.Range("E3").Value = 0.88410487029006224
$ws.Range("C4").Value = 0.95319728792626046
$ws.Range("E4").Value = 0.70119587063581057
$ws.Range("BF4").Value = 0.75866845718531195
$ws.Range("F5").Value = 0.88990852871454607
$ws.Range("G6").Value = 0.72508021818818769
$ws.Range("E7").Value = 0.95860786518984109
$ws.Range("I7").Value = 0.90880477462374609
$ws.Range("AA7").Value = 0.9978391554093724
$ws.Range("F8").Value = 0.59809083114377626
$ws.Range("J8").Value = 0.64506991189634633
$ws.Range("H9").Value = 0.77822181096003096
$ws.Range("J9").Value = 0.88871657634996892
$ws.Range("BP9").Value = 0.71502508624280359
$ws.Range("K10").Value = 0.95204485906003811
$ws.Range("M11").Value = 0.95811087852671717
$ws.Range("J12").Value = 0.89808490947631103
$ws.Range("K12").Value = 0.85652217235921346
$ws.Range("L13").Value = 0.69348561802159292
$ws.Range("L14").Value = 0.85061812785767399
$ws.Range("M14").Value = 0.95592496821695438
$ws.Range("O14").Value = 0.91289311795259032
$ws.Range("M15").Value = 0.67345159593253345
$ws.Range("BA15").Value = 0.76385474510994422
$ws.Range("O17").Value = 0.88673264067572877
$ws.Range("P17").Value = 0.93637518422926558
$ws.Range("T18").Value = 0.9462320942800595
$ws.Range("Q19").Value = 0.78051206746800639
$ws.Range("R19").Value = 0.96809207752073179
$ws.Range("S20").Value = 0.55989210233055364
$ws.Range("V20").Value = 0.7321451363306748
$ws.Range("AE20").Value = 0.94519202898975874
$ws.Range("S21").Value = 0.75167115047312194
$ws.Range("V21").Value = 0.74336452735531044
$ws.Range("W21").Value = 0.976339437079082
$ws.Range("X22").Value = 0.60111531362766069
$ws.Range("V23").Value = 0.88055301619060078
$ws.Range("W24").Value = 0.99933966147390763
$ws.Range("Z25").Value = 0.91274197520936762
$ws.Range("AA25").Value = 0.70903048140348357
$ws.Range("X26").Value = 0.71307882498124875
$ws.Range("Z27").Value = 0.99962501253546399
$ws.Range("AJ27").Value = 0.73652806138303939
$ws.Range("AC28").Value = 0.99836393358033371
$ws.Range("AA29").Value = 0.80114231144894021
$ws.Range("AD29").Value = 0.88928207728569564
$ws.Range("Z30").Value = 0.5960315650879997
$ws.Range("AB30").Value = 0.71022685495811011
$ws.Range("AV30").Value = 0.57587347706475978
$ws.Range("AC31").Value = 0.74241631071169545
$ws.Range("AG31").Value = 0.9554269559455062
$ws.Range("AH32").Value = 0.86246892206561032
$ws.Range("AF33").Value = 0.78344930929703516
$ws.Range("AI33").Value = 0.58224932830541598
$ws.Range("X34").Value = 0.63318864715424017
$ws.Range("AG34").Value = 0.71793801311217997
$ws.Range("AJ35").Value = 0.62629954933422982
$ws.Range("AK35").Value = 0.84598173655468212
$ws.Range("AL35").Value = 0.72221118360089476
$ws.Range("N36").Value = 0.86354617769924846
$ws.Range("AH36").Value = 0.99892886637469824
$ws.Range("AJ37").Value = 0.90622218195658344
$ws.Range("AM37").Value = 0.95555089521551206
$ws.Range("AK38").Value = 0.94452025789172223
$ws.Range("AM38").Value = 0.78202701882515413
$ws.Range("Z39").Value = 0.9236958285543071
$ws.Range("AE40").Value = 0.70395544100031415
$ws.Range("AM40").Value = 0.98462175441094368
$ws.Range("O41").Value = 0.98316629690881729
$ws.Range("AM41").Value = 0.97608817952510518
$ws.Range("AN41").Value = 0.68729479038197083
$ws.Range("AP41").Value = 0.7756712587518424
$ws.Range("AN42").Value = 0.91555136414989835
$ws.Range("AR42").Value = 0.98941329982710702
$ws.Range("AP43").Value = 0.91237206740244048
$ws.Range("AS43").Value = 0.89424305019911199
$ws.Range("AQ44").Value = 0.90185101563572556
$ws.Range("AS44").Value = 0.97843315526193475
$ws.Range("BH45").Value = 0.75594615199753146
$ws.Range("AR46").Value = 0.61405840224338626
$ws.Range("AS46").Value = 0.98718353531670444
$ws.Range("F47").Value = 0.73009917492689203
$ws.Range("AS47").Value = 0.83977125870401792
$ws.Range("AT47").Value = 0.52892488602191789
$ws.Range("AW47").Value = 0.99775919168180194
$ws.Range("N48").Value = 0.69739982736733519
$ws.Range("AT48").Value = 0.6201475275244499
$ws.Range("AV49").Value = 0.96319703018993441
$ws.Range("AX49").Value = 0.72592511587065978
$ws.Range("AV50").Value = 0.94565484141931433
$ws.Range("D51").Value = 0.86687123574640979
$ws.Range("O51").Value = 0.72458490947725784
$ws.Range("AW51").Value = 0.78804955250767228
$ws.Range("AX51").Value = 0.82102743618981022
$ws.Range("BA51").Value = 0.71527138578018434
$ws.Range("AM52").Value = 0.63641326174378721
$ws.Range("AX52").Value = 0.98602771212051876
$ws.Range("AZ53").Value = 0.97942751429402275
$ws.Range("BE53").Value = 0.59003898326872806
$ws.Range("O54").Value = 0.63716363000389364
$ws.Range("AH54").Value = 0.92509447414682122
$ws.Range("BA54").Value = 0.56669196102283159
$ws.Range("C55").Value = 0.86285485999774036
$ws.Range("BB55").Value = 0.76030581760927185
$ws.Range("BE55").Value = 0.95831025333929898
$ws.Range("BD57").Value = 0.81434779971520532
$ws.Range("BF57").Value = 0.97545193767523863
$ws.Range("BD58").Value = 0.85600693355588042
$ws.Range("BG58").Value = 0.81647776248355441
$ws.Range("S59").Value = 0.78016833321346146
$ws.Range("BE59").Value = 0.68413357508290218
$ws.Range("BH59").Value = 0.91585428444991812
$ws.Range("BF60").Value = 0.70477405033494134
$ws.Range("P61").Value = 0.92596348873620737
$ws.Range("BH61").Value = 0.82341757609480792
$ws.Range("W62").Value = 0.92659319305304511
$ws.Range("BI62").Value = 0.92369497264911726
$ws.Range("BK62").Value = 0.68536990093188288
$ws.Range("BL62").Value = 0.90496792216656474
$ws.Range("BI63").Value = 0.94700703976134093
$ws.Range("BM63").Value = 0.75943594956090532
$ws.Range("BK64").Value = 0.8685827875139579
$ws.Range("BM64").Value = 0.66491208101333954
$ws.Range("BN65").Value = 0.92467497501009599
$ws.Range("G66").Value = 0.78285297106500695
$ws.Range("BL66").Value = 0.9096422378911353
$ws.Range("BM67").Value = 0.93285700877550037
$ws.Range("BP67").Value = 0.71722221782593309
